# Scheduled runner: refresh computed profit figures on the Halicarnassus Profits workbook.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) tracks crafted-item profit
# columns H:N; this pass updates the per-row cost/profit figures that changed since
# the last market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4: Root Rush
$ws.Cells.Item(4, 8).Value = 119.75  # H4: 143.5 -> 119.75
$ws.Cells.Item(4, 9).Value = 119.75  # I4: 143.5 -> 119.75
$ws.Cells.Item(4, 11).Value = 119.75  # K4: 143.5 -> 119.75
$ws.Cells.Item(4, 13).Value = -5.75  # M4: -29.5 -> -5.75

# Row 17: One for the Road
$ws.Cells.Item(17, 8).Value = 1749.7576  # H17: 1798.4839 -> 1749.7576
$ws.Cells.Item(17, 9).Value = 699  # I17: 699.9091 -> 699
$ws.Cells.Item(17, 10).Value = 2350.1904  # J17: 2402.7 -> 2350.1904
$ws.Cells.Item(17, 11).Value = 2097  # K17: 2099.7273 -> 2097
$ws.Cells.Item(17, 12).Value = 7050.5712  # L17: 7208.099999999999 -> 7050.5712
$ws.Cells.Item(17, 13).Value = -1929  # M17: -1931.7273 -> -1929
$ws.Cells.Item(17, 14).Value = -7386.5712  # N17: -7544.099999999999 -> -7386.5712

# Row 19: Unbreak My Heart
$ws.Cells.Item(19, 8).Value = 622.2857  # H19: 512.6 -> 622.2857
$ws.Cells.Item(19, 9).Value = 372  # I19: 331.6 -> 372
$ws.Cells.Item(19, 10).Value = 956  # J19: 693.6 -> 956
$ws.Cells.Item(19, 11).Value = 372  # K19: 331.6 -> 372
$ws.Cells.Item(19, 12).Value = 956  # L19: 693.6 -> 956
$ws.Cells.Item(19, 13).Value = -197  # M19: -156.6 -> -197
$ws.Cells.Item(19, 14).Value = -1306  # N19: -1043.6 -> -1306

# Row 106: Making Your Mark
$ws.Cells.Item(106, 8).Value = 1594.5555  # H106: 1678.7142 -> 1594.5555
$ws.Cells.Item(106, 9).Value = 1594.5555  # I106: 1678.7142 -> 1594.5555
$ws.Cells.Item(106, 11).Value = 1594.5555  # K106: 1678.7142 -> 1594.5555
$ws.Cells.Item(106, 13).Value = -963.5554999999999  # M106: -1047.7142 -> -963.5554999999999

# Row 136: I Like Big Brush and I Cannot Lie
$ws.Cells.Item(136, 8).Value = 173333.33  # H136: 245000 -> 173333.33
$ws.Cells.Item(136, 10).Value = 173333.33  # J136: 245000 -> 173333.33
$ws.Cells.Item(136, 12).Value = 173333.33  # L136: 245000 -> 173333.33
$ws.Cells.Item(136, 14).Value = -183533.33  # N136: -255200 -> -183533.33

# Row 141: Remedy for Reason
$ws.Cells.Item(141, 8).Value = 2560  # H141: 2707.1875 -> 2560
$ws.Cells.Item(141, 9).Value = 1602.2142  # I141: 1709.6923 -> 1602.2142
$ws.Cells.Item(141, 11).Value = 4806.642599999999  # K141: 5129.0769 -> 4806.642599999999
$ws.Cells.Item(141, 13).Value = 373.3574000000008  # M141: 50.92309999999998 -> 373.3574000000008

$ws = $wb.Worksheets.Item("ARM")
# Row 88: The Mast Chance
$ws.Cells.Item(88, 8).Value = 678.125  # H88: 676.875 -> 678.125
$ws.Cells.Item(88, 10).Value = 936.25  # J88: 933.75 -> 936.25
$ws.Cells.Item(88, 12).Value = 936.25  # L88: 933.75 -> 936.25
$ws.Cells.Item(88, 14).Value = -1748.25  # N88: -1745.75 -> -1748.25

# Row 91: The Rose and the Riveter (L)
$ws.Cells.Item(91, 8).Value = 678.125  # H91: 676.875 -> 678.125
$ws.Cells.Item(91, 10).Value = 936.25  # J91: 933.75 -> 936.25
$ws.Cells.Item(91, 12).Value = 936.25  # L91: 933.75 -> 936.25
$ws.Cells.Item(91, 14).Value = -3744.25  # N91: -3741.75 -> -3744.25

# Row 97: Ore for Me
$ws.Cells.Item(97, 8).Value = 963.7778  # H97: 980 -> 963.7778
$ws.Cells.Item(97, 9).Value = 835.1539  # I97: 857.61536 -> 835.1539
$ws.Cells.Item(97, 11).Value = 835.1539  # K97: 857.61536 -> 835.1539
$ws.Cells.Item(97, 13).Value = -339.1539  # M97: -361.61536 -> -339.1539

$ws = $wb.Worksheets.Item("BSM")
# Row 29: Powderpost Derby
$ws.Cells.Item(29, 8).Value = 782.2222  # H29: 1003.7 -> 782.2222
$ws.Cells.Item(29, 9).Value = 798.5714  # I29: 1073.375 -> 798.5714
$ws.Cells.Item(29, 11).Value = 798.5714  # K29: 1073.375 -> 798.5714
$ws.Cells.Item(29, 13).Value = -509.5714  # M29: -784.375 -> -509.5714

# Row 36: I Saw What You Did There
$ws.Cells.Item(36, 8).Value = 815  # H36: 991.53845 -> 815
$ws.Cells.Item(36, 10).Value = 0  # J36: 1962.5 -> 0
$ws.Cells.Item(36, 12).Value = 0  # L36: 1962.5 -> 0
$ws.Cells.Item(36, 14).ClearContents()  # N36: cell removed in update

# Row 61: I Maul Right
$ws.Cells.Item(61, 8).Value = 45000  # H61: 44053.5 -> 45000
$ws.Cells.Item(61, 10).Value = 45000  # J61: 44053.5 -> 45000
$ws.Cells.Item(61, 12).Value = 45000  # L61: 44053.5 -> 45000
$ws.Cells.Item(61, 14).Value = -45626  # N61: -44679.5 -> -45626

# Row 94: High Steal
$ws.Cells.Item(94, 8).Value = 385.46667  # H94: 403.5 -> 385.46667
$ws.Cells.Item(94, 9).Value = 377.2857  # I94: 394.72726 -> 377.2857
$ws.Cells.Item(94, 11).Value = 377.2857  # K94: 394.72726 -> 377.2857
$ws.Cells.Item(94, 13).Value = 73.71429999999998  # M94: 56.27274 -> 73.71429999999998

# Row 97: File under Dull
$ws.Cells.Item(97, 8).Value = 1787.25  # H97: 1837.25 -> 1787.25
$ws.Cells.Item(97, 9).Value = 1549.6666  # I97: 1616.3334 -> 1549.6666
$ws.Cells.Item(97, 11).Value = 1549.6666  # K97: 1616.3334 -> 1549.6666
$ws.Cells.Item(97, 13).Value = -558.6666  # M97: -625.3334 -> -558.6666

# Row 103: The Bigger the Blade
$ws.Cells.Item(103, 8).Value = 17999.666  # H103: 0 -> 17999.666
$ws.Cells.Item(103, 10).Value = 17999.666  # J103: 0 -> 17999.666
$ws.Cells.Item(103, 12).Value = 17999.666  # L103: 0 -> 17999.666
$ws.Cells.Item(103, 14).Value = -20343.666  # N103: 0 -> -20343.666

$ws = $wb.Worksheets.Item("CRP")
# Row 33: Tools for the Tools
$ws.Cells.Item(33, 8).Value = 1712.8  # H33: 1613.5834 -> 1712.8
$ws.Cells.Item(33, 9).Value = 904.8333  # I33: 989.8570999999999 -> 904.8333
$ws.Cells.Item(33, 10).Value = 2924.75  # J33: 2486.8 -> 2924.75
$ws.Cells.Item(33, 11).Value = 904.8333  # K33: 989.8570999999999 -> 904.8333
$ws.Cells.Item(33, 12).Value = 2924.75  # L33: 2486.8 -> 2924.75
$ws.Cells.Item(33, 13).Value = -525.8333  # M33: -610.8570999999999 -> -525.8333
$ws.Cells.Item(33, 14).Value = -3682.75  # N33: -3244.8 -> -3682.75

# Row 50: The Arsenal of Theocracy
$ws.Cells.Item(50, 8).Value = 45000  # H50: 44722 -> 45000
$ws.Cells.Item(50, 10).Value = 45000  # J50: 44722 -> 45000
$ws.Cells.Item(50, 12).Value = 45000  # L50: 44722 -> 45000
$ws.Cells.Item(50, 14).Value = -46250  # N50: -45972 -> -46250

# Row 51: Greenstone for Greenhorns
$ws.Cells.Item(51, 8).Value = 33838.75  # H51: 31171 -> 33838.75
$ws.Cells.Item(51, 10).Value = 44755  # J51: 38691.25 -> 44755
$ws.Cells.Item(51, 12).Value = 44755  # L51: 38691.25 -> 44755
$ws.Cells.Item(51, 14).Value = -46227  # N51: -40163.25 -> -46227

# Row 60: Bowing to Greater Power
$ws.Cells.Item(60, 8).Value = 2941.111  # H60: 3669.1428 -> 2941.111
$ws.Cells.Item(60, 9).Value = 2941.111  # I60: 3669.1428 -> 2941.111
$ws.Cells.Item(60, 11).Value = 2941.111  # K60: 3669.1428 -> 2941.111
$ws.Cells.Item(60, 13).Value = -2430.111  # M60: -3158.1428 -> -2430.111

# Row 61: Incant Now, Think Later
$ws.Cells.Item(61, 8).Value = 33838.75  # H61: 31171 -> 33838.75
$ws.Cells.Item(61, 10).Value = 44755  # J61: 38691.25 -> 44755
$ws.Cells.Item(61, 12).Value = 44755  # L61: 38691.25 -> 44755
$ws.Cells.Item(61, 14).Value = -45451  # N61: -39387.25 -> -45451

# Row 99: O Pine
$ws.Cells.Item(99, 8).Value = 3504.1667  # H99: 3605 -> 3504.1667
$ws.Cells.Item(99, 9).Value = 2337.3333  # I99: 2006 -> 2337.3333
$ws.Cells.Item(99, 11).Value = 2337.3333  # K99: 2006 -> 2337.3333
$ws.Cells.Item(99, 13).Value = -839.3332999999998  # M99: -508 -> -839.3332999999998

# Row 109: Playing the Market
$ws.Cells.Item(109, 8).Value = 0  # H109: 15000 -> 0
$ws.Cells.Item(109, 9).Value = 0  # I109: 15000 -> 0
$ws.Cells.Item(109, 11).Value = 0  # K109: 15000 -> 0
$ws.Cells.Item(109, 13).ClearContents()  # M109: cell removed in update

# Row 126: A Better Conductor
$ws.Cells.Item(126, 8).Value = 3504.1667  # H126: 3605 -> 3504.1667
$ws.Cells.Item(126, 9).Value = 2337.3333  # I126: 2006 -> 2337.3333
$ws.Cells.Item(126, 11).Value = 7011.999899999999  # K126: 6018 -> 7011.999899999999
$ws.Cells.Item(126, 13).Value = -4541.999899999999  # M126: -3548 -> -4541.999899999999

# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 1397.0333  # H132: 1507.6296 -> 1397.0333
$ws.Cells.Item(132, 9).Value = 1272.8276  # I132: 1373.3462 -> 1272.8276
$ws.Cells.Item(132, 11).Value = 3818.4828  # K132: 4120.0386 -> 3818.4828
$ws.Cells.Item(132, 13).Value = -1288.4828  # M132: -1590.0386 -> -1288.4828

$ws = $wb.Worksheets.Item("CUL")
# Row 57: The Egg Files
$ws.Cells.Item(57, 8).Value = 1468.75  # H57: 2095 -> 1468.75
$ws.Cells.Item(57, 9).Value = 1329.5454  # I57: 1408.3334 -> 1329.5454
$ws.Cells.Item(57, 10).Value = 3000  # J57: 3125 -> 3000
$ws.Cells.Item(57, 11).Value = 3988.6362  # K57: 4225.0002 -> 3988.6362
$ws.Cells.Item(57, 12).Value = 9000  # L57: 9375 -> 9000
$ws.Cells.Item(57, 13).Value = -3429.6362  # M57: -3666.0002 -> -3429.6362
$ws.Cells.Item(57, 14).Value = -10118  # N57: -10493 -> -10118

# Row 137: Creative Chocolate
$ws.Cells.Item(137, 8).Value = 5775  # H137: 5800 -> 5775
$ws.Cells.Item(137, 10).Value = 5775  # J137: 5800 -> 5775
$ws.Cells.Item(137, 12).Value = 17325  # L137: 17400 -> 17325
$ws.Cells.Item(137, 14).Value = -27525  # N137: -27600 -> -27525

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 1906.3077  # H102: 1820.1428 -> 1906.3077
$ws.Cells.Item(102, 9).Value = 1607.4546  # I102: 1531.8334 -> 1607.4546
$ws.Cells.Item(102, 11).Value = 1607.4546  # K102: 1531.8334 -> 1607.4546
$ws.Cells.Item(102, 13).Value = 14.54539999999997  # M102: 90.16660000000002 -> 14.54539999999997

# Row 122: Awarding Academic Excellence
$ws.Cells.Item(122, 8).Value = 4008  # H122: 2938.6667 -> 4008
$ws.Cells.Item(122, 9).Value = 0  # I122: 800 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 2400 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # M122: cell removed in update

$ws = $wb.Worksheets.Item("LTW")
# Row 34: Breeches Served Cold
$ws.Cells.Item(34, 8).Value = 3572.6  # H34: 2307 -> 3572.6
$ws.Cells.Item(34, 9).Value = 3965.75  # I34: 1960.5 -> 3965.75
$ws.Cells.Item(34, 10).Value = 2000  # J34: 3000 -> 2000
$ws.Cells.Item(34, 11).Value = 3965.75  # K34: 1960.5 -> 3965.75
$ws.Cells.Item(34, 12).Value = 2000  # L34: 3000 -> 2000
$ws.Cells.Item(34, 13).Value = -3793.75  # M34: -1788.5 -> -3793.75
$ws.Cells.Item(34, 14).Value = -2344  # N34: -3344 -> -2344

# Row 40: Best Served Toad
$ws.Cells.Item(40, 8).Value = 8502  # H40: 9377.25 -> 8502
$ws.Cells.Item(40, 9).Value = 8751.25  # I40: 9377.25 -> 8751.25
$ws.Cells.Item(40, 10).Value = 7505  # J40: 0 -> 7505
$ws.Cells.Item(40, 11).Value = 8751.25  # K40: 9377.25 -> 8751.25
$ws.Cells.Item(40, 12).Value = 7505  # L40: 0 -> 7505
$ws.Cells.Item(40, 13).Value = -8615.25  # M40: -9241.25 -> -8615.25
$ws.Cells.Item(40, 14).Value = -7777  # N40: -9241.25 -> -7777

# Row 61: Spelling Me Softly
$ws.Cells.Item(61, 8).Value = 200006400  # H61: 8000 -> 200006400
$ws.Cells.Item(61, 9).Value = 1000000000  # I61: 0 -> 1000000000
$ws.Cells.Item(61, 11).Value = 1000000000  # K61: 0 -> 1000000000
$ws.Cells.Item(61, 13).Value = -999999798  # M61: 0 -> -999999798

# Row 82: Trainin' the Neck
$ws.Cells.Item(82, 8).Value = 3635.9092  # H82: 3153 -> 3635.9092
$ws.Cells.Item(82, 9).Value = 1165.8334  # I82: 998.625 -> 1165.8334
$ws.Cells.Item(82, 11).Value = 1165.8334  # K82: 998.625 -> 1165.8334
$ws.Cells.Item(82, 13).Value = -804.8334  # M82: -637.625 -> -804.8334

# Row 85: Training Is Only Skintight (L)
$ws.Cells.Item(85, 8).Value = 3635.9092  # H85: 3153 -> 3635.9092
$ws.Cells.Item(85, 9).Value = 1165.8334  # I85: 998.625 -> 1165.8334
$ws.Cells.Item(85, 11).Value = 1165.8334  # K85: 998.625 -> 1165.8334
$ws.Cells.Item(85, 13).Value = 82.16660000000002  # M85: 249.375 -> 82.16660000000002

# Row 113: Peace in Rest
$ws.Cells.Item(113, 8).Value = 200006400  # H113: 8000 -> 200006400
$ws.Cells.Item(113, 9).Value = 1000000000  # I113: 0 -> 1000000000
$ws.Cells.Item(113, 11).Value = 1000000000  # K113: 0 -> 1000000000
$ws.Cells.Item(113, 13).Value = -999997830  # M113: 0 -> -999997830

$ws = $wb.Worksheets.Item("WVR")
# Row 18: Welcome to the Cotton Club
$ws.Cells.Item(18, 8).Value = 17666.334  # H18: 16749.25 -> 17666.334
$ws.Cells.Item(18, 9).Value = 15999  # I18: 15998.5 -> 15999
$ws.Cells.Item(18, 10).Value = 18500  # J18: 17500 -> 18500
$ws.Cells.Item(18, 11).Value = 15999  # K18: 15998.5 -> 15999
$ws.Cells.Item(18, 12).Value = 18500  # L18: 17500 -> 18500
$ws.Cells.Item(18, 13).Value = -15826  # M18: -15825.5 -> -15826
$ws.Cells.Item(18, 14).Value = -18846  # N18: -17846 -> -18846

# Row 126: A Polished Purchase
$ws.Cells.Item(126, 8).Value = 5460.643  # H126: 3490.261 -> 5460.643
$ws.Cells.Item(126, 9).Value = 2959.8  # I126: 1417.8462 -> 2959.8
$ws.Cells.Item(126, 10).Value = 6850  # J126: 6184.4 -> 6850
$ws.Cells.Item(126, 11).Value = 8879.400000000001  # K126: 4253.5386 -> 8879.400000000001
$ws.Cells.Item(126, 12).Value = 20550  # L126: 18553.2 -> 20550
$ws.Cells.Item(126, 13).Value = -6409.400000000001  # M126: -1783.5386 -> -6409.400000000001
$ws.Cells.Item(126, 14).Value = -25490  # N126: -23493.2 -> -25490

# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 2963.5  # H132: 4000 -> 2963.5
$ws.Cells.Item(132, 9).Value = 2963.5  # I132: 4000 -> 2963.5
$ws.Cells.Item(132, 11).Value = 8890.5  # K132: 12000 -> 8890.5
$ws.Cells.Item(132, 13).Value = -6360.5  # M132: -9470 -> -6360.5
